$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: copy style from H1 (bold, bordered) to I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF)
$values = @{
    2  = @(9, 9)
    3  = @(8, 8)
    4  = @(4, 4)
    5  = @(8, 8)
    6  = @(5, 6)
    7  = @(5, 6)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(6, 7)
    11 = @(7, 7)
    12 = @(5, 6)
    13 = @(8, 8)
    14 = @(8, 9)
    15 = @(7, 7)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
